$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474990868838/1749474758173/gFRHQATQ91-camera-video-segment-driver-1749474760673.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSCFJKRBTG%2F20250610%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250610T180006Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEOj%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCICFb%2BhOfu4eU2Xric4NOKv%2BHQ1hxJaYhOjrUUid0DcGUAiAPca94rhTMKu%2FhUOWpsFQbz5Pmutplz81F5vH7rxUM7irmAwjB%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM%2FdA6YUe3BT7HzxIQKroDm7hzy03IqvkpIQ8AWo9cVbjDsAo%2F8X10lKtiAdXZME00QMNJwafiO%2BPQE9z0qxwtFHY7us%2BoAQEOvsNrVWh8zE842JY6RdbeaYO%2BqmGjGtOrznN96edVDrVCdqGPsgeUMctIW0hPN0PMZdneMms3woTlqCRS1W71y7ojs4YuTGK5NA8UnoSp5BrgPJ%2F336KFUuHJajqRz4PfA8Ezl4f6dQk9OUzR5d0aAVrAv2YQypuAKge1%2B8qElBFUqx1NOXoIRGVa0RDaoNOAyonKZ5%2FgQ0fPS6f4IaiiuXGEum1QO%2B4SewLQ3arV1CrB6jUNMKk%2BJQJqfpoBGu1nJ47KzlECurFn%2F%2B4UAXiH3YlwozSFx1N2c2ooaTXXTbvv9i3BYfTFIHq32mjA7TLwqXoixA6jTKhHsMZUwyaRDAC30KBo9o9XeeF7nYQAbS2F5ojsZE9DTfTnt3TV8ScsiC2sCOuFCaz4hTSljVzAU9hxGWqMQErYOvbus2i4jXOJTKSDa1k60g%2BeH9SL0r7E0N7Hc5z1jnikcrvPpm0HOUTZ0FaZjlyuvc927Q%2F3BwkhbuQawU6i5hxYhqhBbkXbLDDqqaHCBjqmAUoxQINtu4sMNDdNPYerow40oXs4yjzPePeSXILNWaQgsrlbGgI0LEtOAiwoPAze1SAUqLo8WqHbCpeIt4V%2FSigiEW31zJm5c57hdcOKkd2UAAUxDssvcghRyi2JiGTLPEZZ7N4xGk8ixbKBCBnRszckof9zMgJHtpMs%2BLSa7cejFvB3vqrwzQXmAPcjFPY2v5d02OQ%2BLAQXrQRyPrbXoHdIh1fsc8c%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2011%20Jun%202025%2002%3A00%3A06%20GMT&X-Amz-Signature=0f38052bf2c8f592af1b89e8a1b86a34602ac0d2d1353e182690468834296e80'
$ws.Range("L2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474990868838/1749474758173/7Z3ova6Wtr-camera-video-segment-1749474760673.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSCFJKRBTG%2F20250610%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250610T180006Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEOj%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCICFb%2BhOfu4eU2Xric4NOKv%2BHQ1hxJaYhOjrUUid0DcGUAiAPca94rhTMKu%2FhUOWpsFQbz5Pmutplz81F5vH7rxUM7irmAwjB%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM%2FdA6YUe3BT7HzxIQKroDm7hzy03IqvkpIQ8AWo9cVbjDsAo%2F8X10lKtiAdXZME00QMNJwafiO%2BPQE9z0qxwtFHY7us%2BoAQEOvsNrVWh8zE842JY6RdbeaYO%2BqmGjGtOrznN96edVDrVCdqGPsgeUMctIW0hPN0PMZdneMms3woTlqCRS1W71y7ojs4YuTGK5NA8UnoSp5BrgPJ%2F336KFUuHJajqRz4PfA8Ezl4f6dQk9OUzR5d0aAVrAv2YQypuAKge1%2B8qElBFUqx1NOXoIRGVa0RDaoNOAyonKZ5%2FgQ0fPS6f4IaiiuXGEum1QO%2B4SewLQ3arV1CrB6jUNMKk%2BJQJqfpoBGu1nJ47KzlECurFn%2F%2B4UAXiH3YlwozSFx1N2c2ooaTXXTbvv9i3BYfTFIHq32mjA7TLwqXoixA6jTKhHsMZUwyaRDAC30KBo9o9XeeF7nYQAbS2F5ojsZE9DTfTnt3TV8ScsiC2sCOuFCaz4hTSljVzAU9hxGWqMQErYOvbus2i4jXOJTKSDa1k60g%2BeH9SL0r7E0N7Hc5z1jnikcrvPpm0HOUTZ0FaZjlyuvc927Q%2F3BwkhbuQawU6i5hxYhqhBbkXbLDDqqaHCBjqmAUoxQINtu4sMNDdNPYerow40oXs4yjzPePeSXILNWaQgsrlbGgI0LEtOAiwoPAze1SAUqLo8WqHbCpeIt4V%2FSigiEW31zJm5c57hdcOKkd2UAAUxDssvcghRyi2JiGTLPEZZ7N4xGk8ixbKBCBnRszckof9zMgJHtpMs%2BLSa7cejFvB3vqrwzQXmAPcjFPY2v5d02OQ%2BLAQXrQRyPrbXoHdIh1fsc8c%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2011%20Jun%202025%2002%3A00%3A06%20GMT&X-Amz-Signature=92a39506ff3c0ac9dfd89fc1ccafcc246eb7b1610757defd9aad20019c6f7289'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474990868838/1749474688613/JLWgs3viB6-camera-video-segment-driver-1749474691113.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSCFJKRBTG%2F20250610%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250610T180006Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEOj%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCICFb%2BhOfu4eU2Xric4NOKv%2BHQ1hxJaYhOjrUUid0DcGUAiAPca94rhTMKu%2FhUOWpsFQbz5Pmutplz81F5vH7rxUM7irmAwjB%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM%2FdA6YUe3BT7HzxIQKroDm7hzy03IqvkpIQ8AWo9cVbjDsAo%2F8X10lKtiAdXZME00QMNJwafiO%2BPQE9z0qxwtFHY7us%2BoAQEOvsNrVWh8zE842JY6RdbeaYO%2BqmGjGtOrznN96edVDrVCdqGPsgeUMctIW0hPN0PMZdneMms3woTlqCRS1W71y7ojs4YuTGK5NA8UnoSp5BrgPJ%2F336KFUuHJajqRz4PfA8Ezl4f6dQk9OUzR5d0aAVrAv2YQypuAKge1%2B8qElBFUqx1NOXoIRGVa0RDaoNOAyonKZ5%2FgQ0fPS6f4IaiiuXGEum1QO%2B4SewLQ3arV1CrB6jUNMKk%2BJQJqfpoBGu1nJ47KzlECurFn%2F%2B4UAXiH3YlwozSFx1N2c2ooaTXXTbvv9i3BYfTFIHq32mjA7TLwqXoixA6jTKhHsMZUwyaRDAC30KBo9o9XeeF7nYQAbS2F5ojsZE9DTfTnt3TV8ScsiC2sCOuFCaz4hTSljVzAU9hxGWqMQErYOvbus2i4jXOJTKSDa1k60g%2BeH9SL0r7E0N7Hc5z1jnikcrvPpm0HOUTZ0FaZjlyuvc927Q%2F3BwkhbuQawU6i5hxYhqhBbkXbLDDqqaHCBjqmAUoxQINtu4sMNDdNPYerow40oXs4yjzPePeSXILNWaQgsrlbGgI0LEtOAiwoPAze1SAUqLo8WqHbCpeIt4V%2FSigiEW31zJm5c57hdcOKkd2UAAUxDssvcghRyi2JiGTLPEZZ7N4xGk8ixbKBCBnRszckof9zMgJHtpMs%2BLSa7cejFvB3vqrwzQXmAPcjFPY2v5d02OQ%2BLAQXrQRyPrbXoHdIh1fsc8c%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2011%20Jun%202025%2002%3A00%3A06%20GMT&X-Amz-Signature=de05f4189eebe58536a169141033feae3dfd52f4a19143814db5deacdc8cd803'
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474990868838/1749474688613/LOHXs2uAgt-camera-video-segment-1749474691113.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSCFJKRBTG%2F20250610%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250610T180006Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEOj%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCICFb%2BhOfu4eU2Xric4NOKv%2BHQ1hxJaYhOjrUUid0DcGUAiAPca94rhTMKu%2FhUOWpsFQbz5Pmutplz81F5vH7rxUM7irmAwjB%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIM%2FdA6YUe3BT7HzxIQKroDm7hzy03IqvkpIQ8AWo9cVbjDsAo%2F8X10lKtiAdXZME00QMNJwafiO%2BPQE9z0qxwtFHY7us%2BoAQEOvsNrVWh8zE842JY6RdbeaYO%2BqmGjGtOrznN96edVDrVCdqGPsgeUMctIW0hPN0PMZdneMms3woTlqCRS1W71y7ojs4YuTGK5NA8UnoSp5BrgPJ%2F336KFUuHJajqRz4PfA8Ezl4f6dQk9OUzR5d0aAVrAv2YQypuAKge1%2B8qElBFUqx1NOXoIRGVa0RDaoNOAyonKZ5%2FgQ0fPS6f4IaiiuXGEum1QO%2B4SewLQ3arV1CrB6jUNMKk%2BJQJqfpoBGu1nJ47KzlECurFn%2F%2B4UAXiH3YlwozSFx1N2c2ooaTXXTbvv9i3BYfTFIHq32mjA7TLwqXoixA6jTKhHsMZUwyaRDAC30KBo9o9XeeF7nYQAbS2F5ojsZE9DTfTnt3TV8ScsiC2sCOuFCaz4hTSljVzAU9hxGWqMQErYOvbus2i4jXOJTKSDa1k60g%2BeH9SL0r7E0N7Hc5z1jnikcrvPpm0HOUTZ0FaZjlyuvc927Q%2F3BwkhbuQawU6i5hxYhqhBbkXbLDDqqaHCBjqmAUoxQINtu4sMNDdNPYerow40oXs4yjzPePeSXILNWaQgsrlbGgI0LEtOAiwoPAze1SAUqLo8WqHbCpeIt4V%2FSigiEW31zJm5c57hdcOKkd2UAAUxDssvcghRyi2JiGTLPEZZ7N4xGk8ixbKBCBnRszckof9zMgJHtpMs%2BLSa7cejFvB3vqrwzQXmAPcjFPY2v5d02OQ%2BLAQXrQRyPrbXoHdIh1fsc8c%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2011%20Jun%202025%2002%3A00%3A06%20GMT&X-Amz-Signature=5edf398a34f5025ae4f999645e25333d3e48e7e1d2929da6231b6783a72936e6'
